# Weekly update: insert a new data row at the top of the "Chino" price
# series (row 129), pushing the existing rows 129-194 down to 130-195.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 129; Excel shifts rows 129:194 down to 130:195 and
# copies the row-128 formatting (e.g. the date style on column D) onto it.
$ws.Rows(129).Insert()

# Populate the newly inserted row 129 with this week's record.
$ws.Cells.Item(129, 1).Value  = 7
$ws.Cells.Item(129, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(129, 3).Value  = "Ñuble"
$ws.Cells.Item(129, 4).Value  = 44606
$ws.Cells.Item(129, 5).Value  = 16
$ws.Cells.Item(129, 6).Value  = 100112003
$ws.Cells.Item(129, 7).Value  = "Ajo"
$ws.Cells.Item(129, 8).Value  = "Chino"
$ws.Cells.Item(129, 9).Value  = "Primera"
$ws.Cells.Item(129, 10).Value = 60
$ws.Cells.Item(129, 11).Value = 19000
$ws.Cells.Item(129, 12).Value = 20000
$ws.Cells.Item(129, 13).Value = 19500
$ws.Cells.Item(129, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(129, 15).Value = "China"
$ws.Cells.Item(129, 16).Value = 1950
$ws.Cells.Item(129, 17).Value = 10
$ws.Cells.Item(129, 18).Value = "Hortaliza"
